$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02678272241777
$ws.Range("D2").Value = 1.035990295779749
$ws.Range("E2").Value = 1.026956069319493
$ws.Range("F2").Value = 1.044623112239012
$ws.Range("I2").Value = 1.031752871218361
$ws.Range("J2").Value = 1.03194424383132
$ws.Range("K2").Value = 1.038785318876223
$ws.Range("L2").Value = 1.029777203930575
$ws.Range("M2").Value = 1.047393642020182
$ws.Range("N2").Value = 1.033409723023232

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027783391706685
$ws.Range("D3").Value = 1.036782915075394
$ws.Range("E3").Value = 1.027807269147658
$ws.Range("F3").Value = 1.04565123082715
$ws.Range("I3").Value = 1.031914092206249
$ws.Range("J3").Value = 1.032584556328294
$ws.Range("K3").Value = 1.039387492042086
$ws.Range("L3").Value = 1.030435909376184
$ws.Range("M3").Value = 1.048232470522157
$ws.Range("N3").Value = 1.034050944837385

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028431292805147
$ws.Range("D4").Value = 1.037295916674169
$ws.Range("E4").Value = 1.028358769056981
$ws.Range("F4").Value = 1.04631711992115
$ws.Range("I4").Value = 1.032017072926416
$ws.Range("J4").Value = 1.032998708386763
$ws.Range("K4").Value = 1.039776608543154
$ws.Range("L4").Value = 1.030862212483177
$ws.Range("M4").Value = 1.048775270650242
$ws.Range("N4").Value = 1.03446568503928

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028703765969344
$ws.Range("D5").Value = 1.037511611129003
$ws.Range("E5").Value = 1.02859079018141
$ws.Range("F5").Value = 1.046597209025801
$ws.Range("I5").Value = 1.032060044775318
$ws.Range("J5").Value = 1.033172775813054
$ws.Range("K5").Value = 1.039940065470678
$ws.Range("L5").Value = 1.03104144788709
$ws.Range("M5").Value = 1.04900346821648
$ws.Range("N5").Value = 1.034639999661282

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028749521000924
$ws.Range("D6").Value = 1.037547828821709
$ws.Range("E6").Value = 1.028629757521451
$ws.Range("F6").Value = 1.046644245948887
$ws.Range("I6").Value = 1.032067241087118
$ws.Range("J6").Value = 1.033202000010946
$ws.Range("K6").Value = 1.039967503098991
$ws.Range("L6").Value = 1.0310715433085
$ws.Range("M6").Value = 1.049041783824189
$ws.Range("N6").Value = 1.034669265360887

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028434933230097
$ws.Range("D7").Value = 1.037298798683082
$ws.Range("E7").Value = 1.028361868665684
$ws.Range("F7").Value = 1.046320861899868
$ws.Range("I7").Value = 1.032017648381188
$ws.Range("J7").Value = 1.033001034450526
$ws.Range("K7").Value = 1.039778793165218
$ws.Range("L7").Value = 1.030864607367314
$ws.Range("M7").Value = 1.048778319821128
$ws.Range("N7").Value = 1.034468014406321

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027120819928163
$ws.Range("D8").Value = 1.03625813887381
$ws.Range("E8").Value = 1.02724358727689
$ws.Range("F8").Value = 1.044970439220484
$ws.Range("I8").Value = 1.03180763380668
$ws.Range("J8").Value = 1.032160675827858
$ws.Range("K8").Value = 1.038988935326917
$ws.Range("L8").Value = 1.029999800188894
$ws.Range("M8").Value = 1.047677122959295
$ws.Range("N8").Value = 1.033626462378018

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024808278453296
$ws.Range("D9").Value = 1.034425359983043
$ws.Range("E9").Value = 1.025278566465783
$ws.Range("F9").Value = 1.042595669349913
$ws.Range("I9").Value = 1.0314273174679
$ws.Range("J9").Value = 1.030678562616514
$ws.Range("K9").Value = 1.037593083009602
$ws.Range("L9").Value = 1.028476520237185
$ws.Range("M9").Value = 1.04573688026583
$ws.Range("N9").Value = 1.032142244395913

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023268691582084
$ws.Range("D10").Value = 1.03320424048664
$ws.Range("E10").Value = 1.023972332471113
$ws.Range("F10").Value = 1.041015801609706
$ws.Range("I10").Value = 1.031166911372917
$ws.Range("J10").Value = 1.029689657212544
$ws.Range("K10").Value = 1.036659857809752
$ws.Range("L10").Value = 1.027461462178128
$ws.Range("M10").Value = 1.044443574519446
$ws.Range("N10").Value = 1.03115193463284

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02260253837558
$ws.Range("D11").Value = 1.03267566903902
$ws.Range("E11").Value = 1.023407627378774
$ws.Range("F11").Value = 1.040332497661373
$ws.Range("I11").Value = 1.031052530315924
$ws.Range("J11").Value = 1.029261261456927
$ws.Range("K11").Value = 1.036255140351085
$ws.Range("L11").Value = 1.027022049420865
$ws.Range("M11").Value = 1.043883613887384
$ws.Range("N11").Value = 1.030722930506113

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022355174565809
$ws.Range("D12").Value = 1.032479362306036
$ws.Range("E12").Value = 1.023198007342126
$ws.Range("F12").Value = 1.040078807510529
$ws.Range("I12").Value = 1.031009800567883
$ws.Range("J12").Value = 1.029102107770846
$ws.Range("K12").Value = 1.03610471711644
$ws.Range("L12").Value = 1.026858849780024
$ws.Range("M12").Value = 1.043675627734158
$ws.Range("N12").Value = 1.030563550803542

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022408231536705
$ws.Range("D13").Value = 1.032521469472932
$ws.Range("E13").Value = 1.023242965356087
$ws.Range("F13").Value = 1.040133219495071
$ws.Range("I13").Value = 1.031018977258023
$ws.Range("J13").Value = 1.029136248053613
$ws.Range("K13").Value = 1.03613698761463
$ws.Range("L13").Value = 1.026893855839524
$ws.Range("M13").Value = 1.043720241088795
$ws.Range("N13").Value = 1.030597739569426

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022582089670946
$ws.Range("D14").Value = 1.03265944168147
$ws.Range("E14").Value = 1.02339029732959
$ws.Range("F14").Value = 1.040311525114887
$ws.Range("I14").Value = 1.031049003227361
$ws.Range("J14").Value = 1.029248106340059
$ws.Range("K14").Value = 1.036242708215676
$ws.Range("L14").Value = 1.027008558920127
$ws.Range("M14").Value = 1.043866421516303
$ws.Range("N14").Value = 1.03070975670747

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022689219373545
$ws.Range("D15").Value = 1.032744454684167
$ws.Range("E15").Value = 1.023481091546951
$ws.Range("F15").Value = 1.040421400932074
$ws.Range("I15").Value = 1.031067470959316
$ws.Range("J15").Value = 1.02931702217299
$ws.Range("K15").Value = 1.036307833845323
$ws.Range("L15").Value = 1.027079233655797
$ws.Range("M15").Value = 1.043956489205776
$ws.Range("N15").Value = 1.030778770408788

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023312912497616
$ws.Range("D16").Value = 1.03323932393831
$ws.Range("E16").Value = 1.024009829190087
$ws.Range("F16").Value = 1.041061166937448
$ws.Range("I16").Value = 1.031174468294261
$ws.Range("J16").Value = 1.029718084408553
$ws.Range("K16").Value = 1.036686704472184
$ws.Range("L16").Value = 1.027490627017466
$ws.Range("M16").Value = 1.04448073836509
$ws.Range("N16").Value = 1.031180402198729

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02370427228989
$ws.Range("D17").Value = 1.033549791600275
$ws.Range("E17").Value = 1.024341734722878
$ws.Range("F17").Value = 1.041462687100872
$ws.Range("I17").Value = 1.031241150507962
$ws.Range("J17").Value = 1.029969608851361
$ws.Range("K17").Value = 1.036924193239124
$ws.Range("L17").Value = 1.02774871422885
$ws.Range("M17").Value = 1.044809599739464
$ws.Range("N17").Value = 1.031432283835095

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023932593869946
$ws.Range("D18").Value = 1.033730899566243
$ws.Range("E18").Value = 1.024535416653681
$ws.Range("F18").Value = 1.041696963113108
$ws.Range("I18").Value = 1.031279888445225
$ws.Range("J18").Value = 1.030116300123304
$ws.Range("K18").Value = 1.037062656087707
$ws.Range("L18").Value = 1.027899263094074
$ws.Range("M18").Value = 1.045001423720517
$ws.Range("N18").Value = 1.03157918342547

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024010453766131
$ws.Range("D19").Value = 1.033792655628903
$ws.Range("E19").Value = 1.024601471923439
$ws.Range("F19").Value = 1.041776858061923
$ws.Range("I19").Value = 1.031293070495236
$ws.Range("J19").Value = 1.030166314869597
$ws.Range("K19").Value = 1.037109858085856
$ws.Range("L19").Value = 1.02795059820505
$ws.Range("M19").Value = 1.045066831516112
$ws.Range("N19").Value = 1.03162926919844

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02366227811713
$ws.Range("D20").Value = 1.033516479548161
$ws.Range("E20").Value = 1.024306115389653
$ws.Range("F20").Value = 1.041419599910744
$ws.Range("I20").Value = 1.031234012339825
$ws.Range("J20").Value = 1.029942624586865
$ws.Range("K20").Value = 1.036898719177115
$ws.Range("L20").Value = 1.027721022774827
$ws.Range("M20").Value = 1.044774315534124
$ws.Range("N20").Value = 1.031405261249849

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022530890706547
$ws.Range("D21").Value = 1.032618811517477
$ws.Range("E21").Value = 1.023346907957807
$ws.Range("F21").Value = 1.040259015252734
$ws.Range("I21").Value = 1.031040168047336
$ws.Range("J21").Value = 1.029215167642655
$ws.Range("K21").Value = 1.036211578699332
$ws.Range("L21").Value = 1.026974781224327
$ws.Range("M21").Value = 1.043823374792607
$ws.Range("N21").Value = 1.030676771233337

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021819977504146
$ws.Range("D22").Value = 1.032054575860991
$ws.Range("E22").Value = 1.022744606836845
$ws.Range("F22").Value = 1.039530000585918
$ws.Range("I22").Value = 1.030916881643481
$ws.Range("J22").Value = 1.028757622408988
$ws.Range("K22").Value = 1.03577900762411
$ws.Range("L22").Value = 1.026505692994717
$ws.Range("M22").Value = 1.043225528177851
$ws.Range("N22").Value = 1.030218576232957

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022196804704844
$ws.Range("D23").Value = 1.032353671995052
$ws.Range("E23").Value = 1.023063822798912
$ws.Range("F23").Value = 1.039916399293734
$ws.Range("I23").Value = 1.030982371509366
$ws.Range("J23").Value = 1.029000191157304
$ws.Range("K23").Value = 1.036008372547595
$ws.Range("L23").Value = 1.026754355506543
$ws.Range("M23").Value = 1.043542453218985
$ws.Range("N23").Value = 1.030461489456719

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023681253327695
$ws.Range("D24").Value = 1.033531531776809
$ws.Range("E24").Value = 1.024322209964213
$ws.Range("F24").Value = 1.041439068922184
$ws.Range("I24").Value = 1.031237238254848
$ws.Range("J24").Value = 1.029954817672193
$ws.Range("K24").Value = 1.036910229995656
$ws.Range("L24").Value = 1.027733535317093
$ws.Range("M24").Value = 1.044790258933198
$ws.Range("N24").Value = 1.031417471650756

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025405756409111
$ws.Range("D25").Value = 1.034899052654435
$ws.Range("E25").Value = 1.02578590955212
$ws.Range("F25").Value = 1.043209024381091
$ws.Range("I25").Value = 1.031526849678984
$ws.Range("J25").Value = 1.03106187362352
$ws.Range("K25").Value = 1.037954416053965
$ws.Range("L25").Value = 1.028870246463974
$ws.Range("M25").Value = 1.046238449816381
$ws.Range("N25").Value = 1.032526099748516
